$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Corrected CSV deposit-date derived earning values for rows 42-46 (column I)
$ws.Range("I42").Value = 11.111261485697399
$ws.Range("I43").Value = 8.6074164041601797
$ws.Range("I44").Value = 9.2405371407476196
$ws.Range("I45").Value = 5.5965333255953702
$ws.Range("I46").Value = 5.2242878128563497

# Sum row now needs more decimal precision to show the corrected total
$ws.Range("I47").NumberFormat = "0.00000000000000"

# Widen column I (Gross amount header) so the longer formatted total still fits
$ws.Columns.Item(9).ColumnWidth = 18
